$wb = $excel.ActiveWorkbook

# Rename the "uk_spelling" sheet to "other"
$sheetUk = $wb.Worksheets.Item("uk_spelling")
$sheetUk.Name = "other"

$sheetStats = $wb.Worksheets.Item("stats_terms")

# Fix the "burn-in" typo (was stored as "burn-in-") on row 22
$sheetStats.Range("B22").Value = "burn-in"

# Fix the "k-mean" typo (was stored as "kmean") on row 69
$sheetStats.Range("B69").Value = "k-mean"

# Restore view state: scrolled/selected cells on stats_terms sheet
$sheetStats.Activate()
$sheetStats.Application.ActiveWindow.ScrollRow = 16
$sheetStats.Range("B70").Select()

# Restore view state: selected cell on other sheet
$sheetOther = $wb.Worksheets.Item("other")
$sheetOther.Activate()
$sheetOther.Range("H21").Select()

# Re-activate stats_terms sheet (tabSelected=1 in target)
$sheetStats.Activate()
